# Minor fix in TSP.
# Update column C (Fitness) values on Sheet1 for rows 2 through 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(3789, 3789, 3789, 3789, 4513, 4513, 4513, 4513, 4575, 4575, 4722)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
